$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.046.87'
$ws.Range('E2').Value = '  -0.43%  '
$ws.Range('D3').Value = '2.217.11'
$ws.Range('E3').Value = '  -1.31%  '
$ws.Range('E4').Value = '  +0.26%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '241.82'
$ws.Range('E5').Value = '  -1.97%  '
$ws.Range('E6').Value = '  -0.32%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '73.20'
$ws.Range('E7').Value = '  -1.76%  '
$ws.Range('E8').Value = '  +0.12%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.607'
$ws.Range('E9').Value = '  -2.04%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '42.37'
$ws.Range('E10').Value = '  -0.18%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0955'
$ws.Range('E11').Value = '  +0.79%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.06'
$ws.Range('E12').Value = '  -1.65%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.103'
$ws.Range('E13').Value = '  -0.27%  '
$ws.Range('D14').Value = '2.549.20'
$ws.Range('E14').Value = '  -1.27%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.28'
$ws.Range('E15').Value = '  -1.48%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.837'
$ws.Range('E16').Value = '  -2.03%  '
$ws.Range('D17').Value = '2.220.93'
$ws.Range('E17').Value = '  -1.57%  '
$ws.Range('D18').Value = '41.889.38'
$ws.Range('E18').Value = '  -0.55%  '
$ws.Range('E19').Value = '  +8.23%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.22'
$ws.Range('E20').Value = '  +1.49%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '73.05'
$ws.Range('E21').Value = '  +0.92%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.48'
$ws.Range('E22').Value = '  +15.93%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '230.45'
$ws.Range('E23').Value = '  -0.43%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.08'
$ws.Range('E24').Value = '  -6.53%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '11.92'
$ws.Range('E25').Value = '  +3.80%  '
$ws.Range('E26').Value = '  +0.12%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '3.63'
$ws.Range('E27').Value = '  +0.25%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.28'
$ws.Range('E28').Value = '  -1.39%  '
$ws.Range('E29').Value = '  -2.82%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '167.88'
$ws.Range('E30').Value = '  -0.79%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '20.51'
$ws.Range('E31').Value = '  -1.00%  '
$ws.Range('E32').Value = '  +6.97%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0797'
$ws.Range('E34').Value = '  -0.02%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '29.61'
$ws.Range('E35').Value = '  -4.21%  '
$ws.Range('E36').Value = '  -9.47%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.28'
$ws.Range('E37').Value = '  -3.82%  '
$ws.Range('E38').Value = '  -4.27%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '13.59'
$ws.Range('E39').Value = '  -1.41%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '66.06'
$ws.Range('E40').Value = '  +5.07%  '
$ws.Range('E41').Value = '  -2.14%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.66'
$ws.Range('E42').Value = '  -2.25%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.198'
$ws.Range('E43').Value = '  -3.15%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.81'
$ws.Range('E44').Value = '  +1.54%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '105.15'
$ws.Range('E45').Value = '  -1.62%  '
$ws.Range('E46').Value = '  -2.30%  '
$ws.Range('E47').Value = '  +4.74%  '
$ws.Range('E48').Value = '  -0.03%  '
$ws.Range('E49').Value = '  -0.81%  '
$ws.Range('E50').Value = '  -0.10%  '
$ws.Range('D51').Value = '2.425.25'
$ws.Range('E51').Value = '  -1.34%  '
